# Apply dG value corrections to the Primers Alignment sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("E3").Value = -2.400000095367432
$ws.Range("K3").Value = -0.1000000014901161
$ws.Range("R3").Value = -2.5
$ws.Range("X3").Value = -30.70000076293945

# Row 4
$ws.Range("E4").Value = -2.599999904632568
$ws.Range("R4").Value = -2.5
$ws.Range("X4").Value = -26.39999961853027

# Row 5
$ws.Range("E5").Value = -0.1000000014901161
$ws.Range("K5").Value = -1.299999952316284
$ws.Range("R5").Value = -3.099999904632568
$ws.Range("X5").Value = -26.60000038146973

# Row 6
$ws.Range("E6").Value = -0.699999988079071
$ws.Range("R6").Value = -3.099999904632568
$ws.Range("X6").Value = -30

# Row 7
$ws.Range("X7").Value = -59.90000152587891

# Row 8
$ws.Range("E8").Value = -0.4000000059604645
$ws.Range("R8").Value = -1.600000023841858
$ws.Range("X8").Value = -30.5

# Row 9
$ws.Range("E9").Value = -1.5
$ws.Range("K9").Value = -0.699999988079071
$ws.Range("R9").Value = -1.600000023841858
$ws.Range("X9").Value = -18.60000038146973

# Row 10
$ws.Range("E10").Value = 0
$ws.Range("K10").Value = -0.300000011920929
$ws.Range("R10").Value = -2
$ws.Range("X10").Value = -25.10000038146973

# Row 11
$ws.Range("E11").Value = 0
$ws.Range("K11").Value = -3.900000095367432
$ws.Range("R11").Value = -2
$ws.Range("X11").Value = -24.10000038146973
